# Generate Report for Handback
# - Overview sheet: status text flips from "Ready for handoff" to
#   "Handed back: in sync with en-US" for both tracked files.
# - zh-cn / de-de sheets: populate the "Latest Handback File" (F) and
#   "Latest Handback DateTime" (G) columns by mirroring the existing
#   "Latest Handoff File" (A) and "Latest Target File" (D) hyperlinks,
#   and stamp the handback timestamp into column H.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- Per-language detail sheets -------------------------------------------
function Update-LanguageSheet($sheetName, $handbackStamp) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Pull the existing handoff-file / target-file hyperlink addresses so the
    # new handback-file / handback-datetime columns link to the same assets.
    $addrA2 = ""
    $addrD2 = ""
    $addrA3 = ""
    $addrD3 = ""
    $dispA2 = ""
    $dispD2 = ""
    $dispA3 = ""
    $dispD3 = ""

    foreach ($hl in $ws.Hyperlinks) {
        $rangeAddr = $hl.Range.Address()
        if ($rangeAddr -eq "`$A`$2") {
            $addrA2 = $hl.Address
            $dispA2 = $hl.TextToDisplay
        } elseif ($rangeAddr -eq "`$D`$2") {
            $addrD2 = $hl.Address
            $dispD2 = $hl.TextToDisplay
        } elseif ($rangeAddr -eq "`$A`$3") {
            $addrA3 = $hl.Address
            $dispA3 = $hl.TextToDisplay
        } elseif ($rangeAddr -eq "`$D`$3") {
            $addrD3 = $hl.Address
            $dispD3 = $hl.TextToDisplay
        }
    }

    # New "Latest Handback File" (F) / handback target mirror (G) cells.
    $ws.Hyperlinks.Add($ws.Range("F2"), $addrA2, "", "", $dispA2)
    $ws.Hyperlinks.Add($ws.Range("G2"), $addrD2, "", "", $dispD2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $addrA3, "", "", $dispA3)
    $ws.Hyperlinks.Add($ws.Range("G3"), $addrD3, "", "", $dispD3)

    # "Latest Handback DateTime" (H) now reflects the actual handback time.
    $ws.Range("H2").Value = $handbackStamp
    $ws.Range("H3").Value = $handbackStamp
}

Update-LanguageSheet "zh-cn" "2016-03-24 23:13:17"
Update-LanguageSheet "de-de" "2016-03-24 23:13:24"
